# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Siren_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 778.75
$ws.Range("I21").Value = 778.75
$ws.Range("K21").Value = 778.75
$ws.Range("M21").Value = -310.75
$ws.Range("H23").Value = 778.75
$ws.Range("I23").Value = 778.75
$ws.Range("K23").Value = 778.75
$ws.Range("M23").Value = -544.75
$ws.Range("H58").Value = 691.06665
$ws.Range("I58").Value = 597.5714
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1792.7142
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -1642.7142
$ws.Range("N58").Value = -6300
$ws.Range("H115").Value = 922.4706
$ws.Range("I115").Value = 391.16666
$ws.Range("K115").Value = 1173.49998
$ws.Range("M115").Value = 393.5000199999999
$ws.Range("H137").Value = 9478.879000000001
$ws.Range("I137").Value = 14795.353
$ws.Range("J137").Value = 3830.125
$ws.Range("K137").Value = 44386.05899999999
$ws.Range("L137").Value = 11490.375
$ws.Range("M137").Value = -41836.05899999999
$ws.Range("N137").Value = -16590.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3043.0566
$ws.Range("I32").Value = 2720.204
$ws.Range("K32").Value = 2720.204
$ws.Range("M32").Value = -2433.204
$ws.Range("H61").Value = 4149.691
$ws.Range("I61").Value = 4047.25
$ws.Range("K61").Value = 4047.25
$ws.Range("M61").Value = -3835.25
$ws.Range("H74").Value = 2821.9666
$ws.Range("I74").Value = 919.95654
$ws.Range("J74").Value = 9071.429
$ws.Range("K74").Value = 919.95654
$ws.Range("L74").Value = 9071.429
$ws.Range("M74").Value = -45.95654000000002
$ws.Range("N74").Value = -10819.429
$ws.Range("H77").Value = 2821.9666
$ws.Range("I77").Value = 919.95654
$ws.Range("J77").Value = 9071.429
$ws.Range("K77").Value = 4599.7827
$ws.Range("L77").Value = 45357.145
$ws.Range("M77").Value = -231.7826999999997
$ws.Range("N77").Value = -54093.145
$ws.Range("H122").Value = 2320673.5
$ws.Range("I122").Value = 4929.5
$ws.Range("J122").Value = 6290520.5
$ws.Range("K122").Value = 14788.5
$ws.Range("L122").Value = 18871561.5
$ws.Range("M122").Value = -12338.5
$ws.Range("N122").Value = -18876461.5
$ws.Range("H132").Value = 4086.6567
$ws.Range("I132").Value = 4095.0508
$ws.Range("K132").Value = 12285.1524
$ws.Range("M132").Value = -9755.152399999999
$ws.Range("H136").Value = 4149.691
$ws.Range("I136").Value = 4047.25
$ws.Range("K136").Value = 12141.75
$ws.Range("M136").Value = -9591.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 38299.2
$ws.Range("I26").Value = 38299.2
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 38299.2
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -38007.2
$ws.Range("N26").ClearContents()
$ws.Range("H86").Value = 7020
$ws.Range("I86").Value = 9309.888999999999
$ws.Range("J86").Value = 2898.2
$ws.Range("K86").Value = 9309.888999999999
$ws.Range("L86").Value = 2898.2
$ws.Range("M86").Value = -8186.888999999999
$ws.Range("N86").Value = -5144.2
$ws.Range("H87").Value = 95701.336
$ws.Range("H89").Value = 7020
$ws.Range("I89").Value = 9309.888999999999
$ws.Range("J89").Value = 2898.2
$ws.Range("K89").Value = 46549.44499999999
$ws.Range("L89").Value = 14491
$ws.Range("M89").Value = -40933.44499999999
$ws.Range("N89").Value = -25723
$ws.Range("H90").Value = 95701.336
$ws.Range("H139").Value = 75999
$ws.Range("J139").Value = 75999
$ws.Range("L139").Value = 75999
$ws.Range("N139").Value = -86279

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2762.5186
$ws.Range("I31").Value = 2258.1738
$ws.Range("J31").Value = 5662.5
$ws.Range("K31").Value = 2258.1738
$ws.Range("L31").Value = 5662.5
$ws.Range("M31").Value = -1963.1738
$ws.Range("N31").Value = -6252.5
$ws.Range("H34").Value = 2762.5186
$ws.Range("I34").Value = 2258.1738
$ws.Range("J34").Value = 5662.5
$ws.Range("K34").Value = 2258.1738
$ws.Range("L34").Value = 5662.5
$ws.Range("M34").Value = -2056.1738
$ws.Range("N34").Value = -6066.5
$ws.Range("H92").Value = 56333
$ws.Range("J92").Value = 56333
$ws.Range("L92").Value = 56333
$ws.Range("N92").Value = -61325

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 63464976
$ws.Range("I4").Value = 52497264
$ws.Range("K4").Value = 157491792
$ws.Range("M4").Value = -157491680
$ws.Range("H86").Value = 1059.9
$ws.Range("I86").Value = 849.8333
$ws.Range("J86").Value = 1375
$ws.Range("K86").Value = 2549.4999
$ws.Range("L86").Value = 4125
$ws.Range("M86").Value = -1363.4999
$ws.Range("N86").Value = -6497
$ws.Range("H87").Value = 14671.333
$ws.Range("I87").Value = 14507
$ws.Range("J87").Value = 15000
$ws.Range("K87").Value = 43521
$ws.Range("L87").Value = 45000
$ws.Range("M87").Value = -42273
$ws.Range("N87").Value = -47496
$ws.Range("H89").Value = 1059.9
$ws.Range("I89").Value = 849.8333
$ws.Range("J89").Value = 1375
$ws.Range("K89").Value = 7648.4997
$ws.Range("L89").Value = 12375
$ws.Range("M89").Value = -1720.4997
$ws.Range("N89").Value = -24231
$ws.Range("H90").Value = 14671.333
$ws.Range("I90").Value = 14507
$ws.Range("J90").Value = 15000
$ws.Range("K90").Value = 130563
$ws.Range("L90").Value = 135000
$ws.Range("M90").Value = -124323
$ws.Range("N90").Value = -147480
$ws.Range("H114").Value = 1749.5
$ws.Range("I114").Value = 1749.5
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 5248.5
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -1994.5
$ws.Range("N114").ClearContents()
$ws.Range("H131").Value = 1485.4348
$ws.Range("I131").Value = 709.3333
$ws.Range("J131").Value = 2940.625
$ws.Range("K131").Value = 2127.9999
$ws.Range("L131").Value = 8821.875
$ws.Range("M131").Value = 2912.0001
$ws.Range("N131").Value = -18901.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1282.1538
$ws.Range("I2").Value = 1576.9
$ws.Range("J2").Value = 299.66666
$ws.Range("K2").Value = 1576.9
$ws.Range("L2").Value = 299.66666
$ws.Range("M2").Value = -1463.9
$ws.Range("N2").Value = -525.66666
$ws.Range("H27").Value = 14000
$ws.Range("J27").Value = 14000
$ws.Range("L27").Value = 14000
$ws.Range("N27").Value = -14332
$ws.Range("H33").Value = 23019
$ws.Range("J33").Value = 23019
$ws.Range("L33").Value = 23019
$ws.Range("N33").Value = -23523
$ws.Range("H59").Value = 6020
$ws.Range("I59").Value = 6500
$ws.Range("J59").Value = 5300
$ws.Range("K59").Value = 6500
$ws.Range("L59").Value = 5300
$ws.Range("M59").Value = -5917
$ws.Range("N59").Value = -6466
$ws.Range("H82").Value = 250000
$ws.Range("J82").Value = 250000
$ws.Range("L82").Value = 250000
$ws.Range("N82").Value = -250766
$ws.Range("H85").Value = 250000
$ws.Range("J85").Value = 250000
$ws.Range("L85").Value = 250000
$ws.Range("N85").Value = -252652
$ws.Range("H126").Value = 22620.25
$ws.Range("J126").Value = 18789.455
$ws.Range("L126").Value = 56368.36500000001
$ws.Range("N126").Value = -61308.36500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 50000
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 50000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -48877
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 50000
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 150000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -144384
$ws.Range("N83").ClearContents()
$ws.Range("H136").Value = 5818.3
$ws.Range("I136").Value = 2234.9473
$ws.Range("K136").Value = 6704.841899999999
$ws.Range("M136").Value = -4154.841899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 6166.6665
$ws.Range("I33").Value = 7250
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 7250
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -7000
$ws.Range("N33").Value = -4500
$ws.Range("H36").Value = 6166.6665
$ws.Range("I36").Value = 7250
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 7250
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = -7000
$ws.Range("N36").Value = -4500
$ws.Range("H46").Value = 65000
$ws.Range("I46").Value = 65000
$ws.Range("K46").Value = 65000
$ws.Range("M46").Value = -64769
$ws.Range("H81").Value = 10497.647
$ws.Range("I81").Value = 12897.77
$ws.Range("K81").Value = 25795.54
$ws.Range("M81").Value = -24734.54
$ws.Range("H84").Value = 10497.647
$ws.Range("I84").Value = 12897.77
$ws.Range("K84").Value = 128977.7
$ws.Range("M84").Value = -123673.7
$ws.Range("H96").Value = 2314.5
$ws.Range("I96").Value = 2099.75
$ws.Range("K96").Value = 2099.75
$ws.Range("M96").Value = -726.75
$ws.Range("H113").Value = 2436.0625
$ws.Range("I113").Value = 1213.1428
$ws.Range("K113").Value = 3639.4284
$ws.Range("M113").Value = -1469.4284
$ws.Range("H134").Value = 65000
$ws.Range("I134").Value = 65000
$ws.Range("K134").Value = 195000
$ws.Range("M134").Value = -192465

